$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-10 from 2023-09-21 (45190)
# to 2023-09-23 (45192), keeping the existing date formatting.
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45192
}
